$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.835.29'
$ws.Range("E2").Value = '  +3.46%  '
$ws.Range("D3").Value = '2.424.08'
$ws.Range("E3").Value = '  +3.22%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.09'
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.22'
$ws.Range("E6").Value = '  +3.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +3.51%  '
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.70'
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.77'
$ws.Range("E13").Value = '  +4.19%  '
$ws.Range("D14").Value = '2.853.64'
$ws.Range("E14").Value = '  +3.15%  '
$ws.Range("D15").Value = '59.783.91'
$ws.Range("E15").Value = '  +3.47%  '
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = '2.421.87'
$ws.Range("E17").Value = '  +3.91%  '
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '330.43'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.74'
$ws.Range("E23").Value = '  +3.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.171'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.65'
$ws.Range("E25").Value = '  +5.42%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.36'
$ws.Range("E27").Value = '  +3.42%  '
$ws.Range("D28").Value = '0.0₃0776'
$ws.Range("E28").Value = '  +5.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.77'
$ws.Range("E29").Value = '  +0.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.55'
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.15'
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.67'
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("E35").Value = '  +4.76%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.60'
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.58'
$ws.Range("E39").Value = '  +1.16%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.410'
$ws.Range("E40").Value = '  -3.46%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '313.67'
$ws.Range("E41").Value = '  +8.59%  '
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '138.60'
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0968'
$ws.Range("E44").Value = '  +2.06%  '
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.47'
$ws.Range("E46").Value = '  +4.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.581'
$ws.Range("E47").Value = '  +3.03%  '
$ws.Range("B48").Value = 'Polygon'
$ws.Range("C48").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.402'
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0225'
$ws.Range("E49").Value = '  +1.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.61'
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.04'

Write-Host "Applied crypto price/volume updates"
